$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 31.03919159155087
$ws.Cells.Item(2, 3).Value = 11.80897108301087
$ws.Cells.Item(2, 4).Value = 3.638390504451034
$ws.Cells.Item(2, 5).Value = 9.781564457095364
$ws.Cells.Item(2, 6).Value = 57.21429344714861
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 10).Value = 9.924875537058677
$ws.Cells.Item(2, 12).Value = 12.52259699289039
$ws.Cells.Item(2, 14).Value = 22.52100569876992

$ws.Cells.Item(3, 2).Value = 30.75136037658014
$ws.Cells.Item(3, 3).Value = 11.45441418776792
$ws.Cells.Item(3, 4).Value = 3.598486616253938
$ws.Cells.Item(3, 5).Value = 9.791276608609017
$ws.Cells.Item(3, 6).Value = 57.1017890829133
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 10).Value = 9.944898989191147
$ws.Cells.Item(3, 12).Value = 12.53141073717467
$ws.Cells.Item(3, 14).Value = 22.58261734733783

$ws.Cells.Item(4, 2).Value = 30.58223800016076
$ws.Cells.Item(4, 3).Value = 11.23524848936645
$ws.Cells.Item(4, 4).Value = 3.573330323438677
$ws.Cells.Item(4, 5).Value = 9.797639032715338
$ws.Cells.Item(4, 6).Value = 57.04642659920724
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 10).Value = 9.957879266023536
$ws.Cells.Item(4, 12).Value = 12.53898630945135
$ws.Cells.Item(4, 14).Value = 22.62246927595953

$ws.Cells.Item(5, 2).Value = 30.51529977642626
$ws.Cells.Item(5, 3).Value = 11.14572470794387
$ws.Cells.Item(5, 4).Value = 3.56291594669956
$ws.Cells.Item(5, 5).Value = 9.800332367482326
$ws.Cells.Item(5, 6).Value = 57.02731752222397
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 10).Value = 9.963341812148711
$ws.Cells.Item(5, 12).Value = 12.54261745264082
$ws.Cells.Item(5, 14).Value = 22.63921799608947

$ws.Cells.Item(6, 2).Value = 30.50430624710864
$ws.Cells.Item(6, 3).Value = 11.13085129325346
$ws.Cells.Item(6, 4).Value = 3.561176855484029
$ws.Cells.Item(6, 5).Value = 9.800785676260654
$ws.Cells.Item(6, 6).Value = 57.02435295723826
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 10).Value = 9.964259326765804
$ws.Cells.Item(6, 12).Value = 12.54325325405579
$ws.Cells.Item(6, 14).Value = 22.6420298426204

$ws.Cells.Item(7, 2).Value = 30.58132714253229
$ws.Cells.Item(7, 3).Value = 11.23404178532485
$ws.Cells.Item(7, 4).Value = 3.57319052909174
$ws.Cells.Item(7, 5).Value = 9.797674948349954
$ws.Cells.Item(7, 6).Value = 57.04615491090337
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 10).Value = 9.957952234740336
$ws.Cells.Item(7, 12).Value = 12.53903307775089
$ws.Cells.Item(7, 14).Value = 22.62269309482576

$ws.Cells.Item(8, 2).Value = 30.93841035837426
$ws.Cells.Item(8, 3).Value = 11.68711554888616
$ws.Cells.Item(8, 4).Value = 3.624768659502345
$ws.Cells.Item(8, 5).Value = 9.784830507987843
$ws.Cells.Item(8, 6).Value = 57.17265778689392
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 10).Value = 9.931637605054732
$ws.Cells.Item(8, 12).Value = 12.52518683756887
$ws.Cells.Item(8, 14).Value = 22.54182985460929

$ws.Cells.Item(9, 2).Value = 31.69584672181441
$ws.Cells.Item(9, 3).Value = 12.55753325733104
$ws.Cells.Item(9, 4).Value = 3.72064164927088
$ws.Cells.Item(9, 5).Value = 9.762799133512143
$ws.Cells.Item(9, 6).Value = 57.52930836514225
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 10).Value = 9.885452127500228
$ws.Cells.Item(9, 12).Value = 12.5152040514705
$ws.Cells.Item(9, 14).Value = 22.39928959406292

$ws.Cells.Item(10, 2).Value = 32.28274083786071
$ws.Cells.Item(10, 3).Value = 13.17813083140936
$ws.Cells.Item(10, 4).Value = 3.787773898280446
$ws.Cells.Item(10, 5).Value = 9.748522915615609
$ws.Cells.Item(10, 6).Value = 57.85700943273506
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 10).Value = 9.854788709180854
$ws.Cells.Item(10, 12).Value = 12.51833087046732
$ws.Cells.Item(10, 14).Value = 22.30432079342258

$ws.Cells.Item(11, 2).Value = 32.55533601695991
$ws.Cells.Item(11, 3).Value = 13.4548988771846
$ws.Cells.Item(11, 4).Value = 3.817577881972676
$ws.Cells.Item(11, 5).Value = 9.742440155649797
$ws.Cells.Item(11, 6).Value = 58.02017144506028
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 10).Value = 9.841541836948039
$ws.Cells.Item(11, 12).Value = 12.52202057339672
$ws.Cells.Item(11, 14).Value = 22.26323205396104

$ws.Cells.Item(12, 2).Value = 32.65928184564266
$ws.Cells.Item(12, 3).Value = 13.55879991973061
$ws.Cells.Item(12, 4).Value = 3.828756638737781
$ws.Cells.Item(12, 5).Value = 9.740195734021581
$ws.Cells.Item(12, 6).Value = 58.08396225767282
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 10).Value = 9.836626004886464
$ws.Cells.Item(12, 12).Value = 12.52374311024321
$ws.Cells.Item(12, 14).Value = 22.24797657958243

$ws.Cells.Item(13, 2).Value = 32.63686462911244
$ws.Cells.Item(13, 3).Value = 13.53646482316551
$ws.Cells.Item(13, 4).Value = 3.826353894957728
$ws.Cells.Item(13, 5).Value = 9.740676489950474
$ws.Cells.Item(13, 6).Value = 58.0701349648407
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 10).Value = 9.83768025650032
$ws.Cells.Item(13, 12).Value = 12.5233576751887
$ws.Cells.Item(13, 14).Value = 22.25124859932913

$ws.Cells.Item(14, 2).Value = 32.5638737766181
$ws.Cells.Item(14, 3).Value = 13.46346559512065
$ws.Cells.Item(14, 4).Value = 3.818499726286141
$ws.Cells.Item(14, 5).Value = 9.742254324633222
$ws.Cells.Item(14, 6).Value = 58.0253795042188
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 10).Value = 9.841135397669424
$ws.Cells.Item(14, 12).Value = 12.52215577216572
$ws.Cells.Item(14, 14).Value = 22.26197088487606

$ws.Cells.Item(15, 2).Value = 32.51925592667943
$ws.Cells.Item(15, 3).Value = 13.41863058756197
$ws.Cells.Item(15, 4).Value = 3.81367478971517
$ws.Cells.Item(15, 5).Value = 9.743228470149598
$ws.Cells.Item(15, 6).Value = 57.99822591066529
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 10).Value = 9.843264841323496
$ws.Cells.Item(15, 12).Value = 12.52146191508698
$ws.Cells.Item(15, 14).Value = 22.26857817947871

$ws.Cells.Item(16, 2).Value = 32.26503162235838
$ws.Cells.Item(16, 3).Value = 13.15992236424344
$ws.Cells.Item(16, 4).Value = 3.785811204835102
$ws.Cells.Item(16, 5).Value = 9.748928703249195
$ws.Cells.Item(16, 6).Value = 57.84662818772654
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 10).Value = 9.855668509266811
$ws.Cells.Item(16, 12).Value = 12.51813530680827
$ws.Cells.Item(16, 14).Value = 22.30704856551477

$ws.Cells.Item(17, 2).Value = 32.11044926723483
$ws.Cells.Item(17, 3).Value = 12.99971306552999
$ws.Cells.Item(17, 4).Value = 3.768528366144047
$ws.Cells.Item(17, 5).Value = 9.752530882977654
$ws.Cells.Item(17, 6).Value = 57.7572224328881
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 10).Value = 9.86345722797849
$ws.Cells.Item(17, 12).Value = 12.51667479884281
$ws.Cells.Item(17, 14).Value = 22.33119015162918

$ws.Cells.Item(18, 2).Value = 32.02207130839243
$ws.Cells.Item(18, 3).Value = 12.90704874918326
$ws.Cells.Item(18, 4).Value = 3.758518668300754
$ws.Cells.Item(18, 5).Value = 9.7546415111489
$ws.Cells.Item(18, 6).Value = 57.70712620813963
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 10).Value = 9.868003204021582
$ws.Cells.Item(18, 12).Value = 12.51604821508308
$ws.Cells.Item(18, 14).Value = 22.34527463568505

$ws.Cells.Item(19, 2).Value = 31.99224236586517
$ws.Cells.Item(19, 3).Value = 12.87558900770892
$ws.Cells.Item(19, 4).Value = 3.755117747241057
$ws.Cells.Item(19, 5).Value = 9.755362794876008
$ws.Cells.Item(19, 6).Value = 57.69039309604549
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 10).Value = 9.869553762890702
$ws.Cells.Item(19, 12).Value = 12.51587274567
$ws.Cells.Item(19, 14).Value = 22.3500775561393

$ws.Cells.Item(20, 2).Value = 32.12685021492735
$ws.Cells.Item(20, 3).Value = 13.01682185374424
$ws.Cells.Item(20, 4).Value = 3.77037531441435
$ws.Cells.Item(20, 5).Value = 9.752143415777292
$ws.Cells.Item(20, 6).Value = 57.76660256960751
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 10).Value = 9.862621266867293
$ws.Cells.Item(20, 12).Value = 12.51680818452782
$ws.Cells.Item(20, 14).Value = 22.32859965452628

$ws.Cells.Item(21, 2).Value = 32.58529411942742
$ws.Cells.Item(21, 3).Value = 13.48493260495863
$ws.Cells.Item(21, 4).Value = 3.820809613109696
$ws.Cells.Item(21, 5).Value = 9.741789276966589
$ws.Cells.Item(21, 6).Value = 58.03847103884374
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 10).Value = 9.840117816055386
$ws.Cells.Item(21, 12).Value = 12.52249997742588
$ws.Cells.Item(21, 14).Value = 22.25881323997892

$ws.Cells.Item(22, 2).Value = 32.88906928533746
$ws.Cells.Item(22, 3).Value = 13.78555242829372
$ws.Cells.Item(22, 4).Value = 3.853144793930692
$ws.Cells.Item(22, 5).Value = 9.735365985914243
$ws.Cells.Item(22, 6).Value = 58.22782696162041
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 10).Value = 9.825995908763927
$ws.Cells.Item(22, 12).Value = 12.52811569771226
$ws.Cells.Item(22, 14).Value = 22.21497544384427

$ws.Cells.Item(23, 2).Value = 32.72658752428654
$ws.Cells.Item(23, 3).Value = 13.62562537331313
$ws.Cells.Item(23, 4).Value = 3.835944773878585
$ws.Cells.Item(23, 5).Value = 9.738762828522201
$ws.Cells.Item(23, 6).Value = 58.1257037129413
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 10).Value = 9.833479632038465
$ws.Cells.Item(23, 12).Value = 12.52494529305484
$ws.Cells.Item(23, 14).Value = 22.23821035487266

$ws.Cells.Item(24, 2).Value = 32.11943380272872
$ws.Cells.Item(24, 3).Value = 13.00908869990662
$ws.Cells.Item(24, 4).Value = 3.769540538106126
$ws.Cells.Item(24, 5).Value = 9.752318466067136
$ws.Cells.Item(24, 6).Value = 57.76235774552136
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 10).Value = 9.862998992614285
$ws.Cells.Item(24, 12).Value = 12.51674721710557
$ws.Cells.Item(24, 14).Value = 22.32977017923061

$ws.Cells.Item(25, 2).Value = 31.4852667414178
$ws.Cells.Item(25, 3).Value = 12.3248511158275
$ws.Cells.Item(25, 4).Value = 3.695279473714009
$ws.Cells.Item(25, 5).Value = 9.768422744724461
$ws.Cells.Item(25, 6).Value = 57.42124195434435
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 10).Value = 9.897370032262279
$ws.Cells.Item(25, 12).Value = 12.51606651179702
$ws.Cells.Item(25, 14).Value = 22.43613615576073
